# "Fruta / hortaliza, semanal" weekly update:
# A new weekly price-report row is inserted at row 143 of the data table
# (the sheet is ordered with the most recent report first), pushing every
# existing row down by one (old row 143 -> new row 144, ..., old row 157 ->
# new row 158). The worksheet dimension grows from A1:T157 to A1:T158.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 143, shifting rows
# 143-157 down to 144-158 (and all their formatting/styles with them).
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new weekly record.
$ws.Range("A143").Value = 9
$ws.Range("B143").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C143").Value = "Metropolitana"
$ws.Range("D143").Value = 45013
$ws.Range("E143").Value = 13
$ws.Range("F143").Value = "Fruta"
$ws.Range("G143").Value = 100101
$ws.Range("H143").Value = "Berries"
$ws.Range("I143").Value = 100101004
$ws.Range("J143").Value = "Frambuesa"
$ws.Range("K143").Value = "Sin especificar"
$ws.Range("L143").Value = "Primera"
$ws.Range("M143").Value = 500
$ws.Range("N143").Value = 7500
$ws.Range("O143").Value = 8000
$ws.Range("P143").Value = 7720
$ws.Range("Q143").Value = "$/bandeja 2 kilos"
$ws.Range("R143").Value = "Provincia de Curicó"
$ws.Range("S143").Value = 3860
$ws.Range("T143").Value = 2
